$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.80"
$ws.Range("E2").Value = "'-1.06%"
$ws.Range("D3").Value = "'31.53"
$ws.Range("E3").Value = "'-1.82%"
$ws.Range("D4").Value = "'5.147"
$ws.Range("E4").Value = "'-2.79%"
$ws.Range("D5").Value = "'0.07408"
$ws.Range("E5").Value = "'-0.94%"
$ws.Range("D6").Value = "'2.255"
$ws.Range("E6").Value = "'43.48%"
$ws.Range("D7").Value = "'7.936"
$ws.Range("E7").Value = "'1.02%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9263"
$ws.Range("E8").Value = "'0.91%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1728"
$ws.Range("E9").Value = "'1.89%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.07626"
$ws.Range("E10").Value = "'-2.86%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08070"
$ws.Range("E11").Value = "'0.60%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03037"
$ws.Range("E12").Value = "'0.72%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09943"
$ws.Range("E13").Value = "'0.33%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001501"
$ws.Range("E14").Value = "'-1.03%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006082"
$ws.Range("E15").Value = "'-4.87%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.458"
$ws.Range("E16").Value = "'-0.45%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "'3.765"
$ws.Range("E17").Value = "'-1.09%"
$ws.Range("D18").Value = "'2.232"
$ws.Range("D19").Value = "'0.3248"
$ws.Range("E19").Value = "'-2.40%"
$ws.Range("D20").Value = "'0.1323"
$ws.Range("E20").Value = "'-0.83%"
$ws.Range("D21").Value = "'4.656"
$ws.Range("E21").Value = "'3.75%"
$ws.Range("D22").Value = "'0.04660"
$ws.Range("E22").Value = "'1.13%"
$ws.Range("D23").Value = "'0.1583"
$ws.Range("E23").Value = "'-2.31%"
$ws.Range("D24").Value = "'0.001226"
$ws.Range("E24").Value = "'0.52%"
$ws.Range("D25").Value = "'0.004496"
$ws.Range("E25").Value = "'1.08%"
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'-6.88%"
$ws.Range("E27").Value = "'7.16%"
$ws.Range("D39").Value = "'0.01733"
$ws.Range("E39").Value = "'-1.60%"
$ws.Range("D40").Value = "'0.04531"
$ws.Range("E40").Value = "'0.08%"
$ws.Range("D41").Value = "'0.007095"
$ws.Range("E41").Value = "'-1.11%"
$ws.Range("D42").Value = "'0.1347"
$ws.Range("E42").Value = "'-0.08%"
$ws.Range("D43").Value = "'0.002224"
$ws.Range("E43").Value = "'0.29%"
$ws.Range("E44").Value = "'-13.93%"
$ws.Range("D45").Value = "'0.00006273"
$ws.Range("E45").Value = "'0.91%"
$ws.Range("E46").Value = "'-46.12%"
$ws.Range("D47").Value = "'1.928"
$ws.Range("E47").Value = "'2.98%"
